# Applies the commit's changes to the "Session Analysis Results" sheet:
#  1. Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" in
#     every "Recorded By" (column G) cell that has that exact value.
#  2. Update the "Class Statistics" summary: Missing Sessions (L7) and
#     Pending Sessions (L8).
#  3. Update the "Group Statistics" table: Missing (P) / Pending (Q) counts
#     for each of the 12 groups (rows 15-26).
#  4. Re-classify the 12 "B1-N / SURGERY SEMINAR/SLIDE / session 8" rows
#     from "Pending" to "Not Recorded" (text in column I) and restyle
#     columns A-I on those rows to match the existing "Not Recorded" look
#     (pink fill) instead of the "Pending" look (pale yellow fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Recorded By" text swap (column G)
# ---------------------------------------------------------------------------
$gRows = @(2,3,4,5,6,7,8,16,17,22,23,24,26,29,37,38,43,44,45,47,50,58,59,64,65,
           66,68,71,79,80,85,86,87,88,89,90,91,99,100,105,106,107,108,109,110,
           111,119,120,125,126,127,128,129,130,131,139,140,145,146,147,148,
           149,150,151,159,160,165,166,167,168,169,170,171,179,180,185,186,
           187,189,192,200,201,206,207,208,210,213,221,222,227,228,229,231,
           234,242,243)

foreach ($r in $gRows) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}

# ---------------------------------------------------------------------------
# 2) Class Statistics totals
# ---------------------------------------------------------------------------
$ws.Range("L7").Value = 84   # Missing Sessions: 72 -> 84
$ws.Range("L8").Value = 24   # Pending Sessions: 36 -> 24

# ---------------------------------------------------------------------------
# 3) Group Statistics table (rows 15-26): Missing +1, Pending -1
# ---------------------------------------------------------------------------
for ($r = 15; $r -le 26; $r++) {
    $pCell = $ws.Range("P$r")
    $qCell = $ws.Range("Q$r")
    $pCell.Value = $pCell.Value2 + 1
    $qCell.Value = $qCell.Value2 - 1
}

# ---------------------------------------------------------------------------
# 4) "Pending" -> "Not Recorded" rows (session 8 rows for every B1-N group)
# ---------------------------------------------------------------------------
$statusRows = @(21,42,63,84,104,124,144,164,184,205,226,247)

# Row 10 is an existing "Not Recorded" row; copy its A:I formatting
# (pink fill) onto each target row so the style resolves the same way
# Excel would (matching cell style index used elsewhere in the sheet).
$formatSource = $ws.Range("A10:I10")

foreach ($r in $statusRows) {
    $formatSource.Copy()
    $destRow = $ws.Range("A" + $r + ":I" + $r)
    $destRow.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Range("I$r").Value = "Not Recorded"
}
$excel.CutCopyMode = 0
